$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "31.531.86"
$ws.Range("E2").Value = "  +5.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.63"
$ws.Range("E3").Value = "  +4.24%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "221.69"
$ws.Range("E5").Value = "  +2.90%  "
$ws.Range("E6").Value = "  +2.95%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.79"
$ws.Range("E8").Value = "  +3.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.43"
$ws.Range("E9").Value = "  +3.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.269"
$ws.Range("E10").Value = "  +3.32%  "
$ws.Range("E11").Value = "  +5.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0911"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.958.01"
$ws.Range("E13").Value = "  +4.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.706.37"
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.611"
$ws.Range("E15").Value = "  +3.63%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "10.19"
$ws.Range("E16").Value = "  +7.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.17"
$ws.Range("E17").Value = "  +7.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "31.544.39"
$ws.Range("E18").Value = "  +5.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.20"
$ws.Range("E19").Value = "  +4.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "251.13"
$ws.Range("E20").Value = "  +4.20%  "
$ws.Range("E21").Value = "  +2.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.16"
$ws.Range("E23").Value = "  +2.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.25"
$ws.Range("E24").Value = "  +2.44%  "
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.91"
$ws.Range("E26").Value = "  +0.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.02"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("E28").Value = "  +3.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.78"
$ws.Range("E29").Value = "  +2.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("E31").Value = "  +12.16%  "
$ws.Range("E32").Value = "  +1.85%  "
$ws.Range("E33").Value = "  +4.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.39"
$ws.Range("E34").Value = "  +5.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.518.10"
$ws.Range("E35").Value = "  +6.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.73"
$ws.Range("E36").Value = "  +1.78%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.03"
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "82.91"
$ws.Range("E38").Value = "  +7.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.606"
$ws.Range("E39").Value = "  +7.65%  "
$ws.Range("E40").Value = "  +4.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +1.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.31"
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.852"
$ws.Range("E43").Value = "  +2.18%  "
$ws.Range("E44").Value = "  +4.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0504"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("E46").Value = "  +3.52%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.999"
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.18"
$ws.Range("E48").Value = "  +6.49%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.58"
$ws.Range("E49").Value = "  +3.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.847.23"
$ws.Range("E50").Value = "  +3.69%  "
$ws.Range("E51").Value = "  +10.24%  "
